$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72..158 down to 73..159
$ws.Rows("72:72").Insert()

# Fill in the values for the newly inserted row 72
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = "Vega Monumental Concepción"
$ws.Range("C72").Value = "Bíobío"
$ws.Range("D72").Value = 44650
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100108
$ws.Range("H72").Value = "Tropicales y subtropicales"
$ws.Range("I72").Value = 100108005
$ws.Range("J72").Value = "Piña"
$ws.Range("K72").Value = "Sin especificar"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 180
$ws.Range("N72").Value = 14000
$ws.Range("O72").Value = 16000
$ws.Range("P72").Value = 15111
$ws.Range("Q72").Value = "$/caja 14 unidades"
$ws.Range("R72").Value = "Ecuador"
$ws.Range("S72").Value = 1079
$ws.Range("T72").Value = 14
